$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 217.86957
$ws.Range("I33").Value = 218.59091
$ws.Range("J33").Value = 202
$ws.Range("K33").Value = 218.59091
$ws.Range("L33").Value = 202
$ws.Range("M33").Value = 10.40908999999999
$ws.Range("N33").Value = -660
$ws.Range("H125").Value = 67438.2
$ws.Range("I125").Value = 111713.664
$ws.Range("J125").Value = 1025
$ws.Range("K125").Value = 1005422.976
$ws.Range("L125").Value = 9225
$ws.Range("M125").Value = -1002962.976
$ws.Range("N125").Value = -14145
$ws.Range("H127").Value = 1175.2609
$ws.Range("I127").Value = 792.4286
$ws.Range("J127").Value = 1218.4839
$ws.Range("K127").Value = 2377.2858
$ws.Range("L127").Value = 3655.4517
$ws.Range("M127").Value = 2582.7142
$ws.Range("N127").Value = -13575.4517
$ws.Range("H131").Value = 4804.222
$ws.Range("I131").Value = 1061.8182
$ws.Range("J131").Value = 6015
$ws.Range("K131").Value = 3185.4546
$ws.Range("L131").Value = 18045
$ws.Range("M131").Value = 1854.5454
$ws.Range("N131").Value = -28125
$ws.Range("H138").Value = 2384908.2
$ws.Range("I138").Value = 13336383
$ws.Range("J138").Value = 4152.681
$ws.Range("K138").Value = 40009149
$ws.Range("L138").Value = 12458.043
$ws.Range("M138").Value = -40004009
$ws.Range("N138").Value = -22738.043
$ws.Range("H141").Value = 9088.044
$ws.Range("I141").Value = 4955.6816
$ws.Range("J141").Value = 100000
$ws.Range("K141").Value = 14867.0448
$ws.Range("L141").Value = 300000
$ws.Range("M141").Value = -9687.0448
$ws.Range("N141").Value = -310360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 916.6486
$ws.Range("I2").Value = 611.8889
$ws.Range("J2").Value = 1739.5
$ws.Range("K2").Value = 611.8889
$ws.Range("L2").Value = 1739.5
$ws.Range("M2").Value = -498.8889
$ws.Range("N2").Value = -1965.5
$ws.Range("H34").Value = 19999.834
$ws.Range("I34").Value = 19999
$ws.Range("K34").Value = 19999
$ws.Range("M34").Value = -19728
$ws.Range("H97").Value = 687.1429000000001
$ws.Range("I97").Value = 687.1429000000001
$ws.Range("K97").Value = 687.1429000000001
$ws.Range("M97").Value = -191.1429000000001
$ws.Range("H116").Value = 916.6486
$ws.Range("I116").Value = 611.8889
$ws.Range("J116").Value = 1739.5
$ws.Range("K116").Value = 611.8889
$ws.Range("L116").Value = 1739.5
$ws.Range("M116").Value = 1682.1111
$ws.Range("N116").Value = -6327.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 916.6486
$ws.Range("I3").Value = 611.8889
$ws.Range("J3").Value = 1739.5
$ws.Range("K3").Value = 611.8889
$ws.Range("L3").Value = 1739.5
$ws.Range("M3").Value = -497.8889
$ws.Range("N3").Value = -1967.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 14086820
$ws.Range("I31").Value = 24392010
$ws.Range("J31").Value = 3060.2
$ws.Range("K31").Value = 24392010
$ws.Range("L31").Value = 3060.2
$ws.Range("M31").Value = -24391715
$ws.Range("N31").Value = -3650.2
$ws.Range("H34").Value = 14086820
$ws.Range("I34").Value = 24392010
$ws.Range("J34").Value = 3060.2
$ws.Range("K34").Value = 24392010
$ws.Range("L34").Value = 3060.2
$ws.Range("M34").Value = -24391808
$ws.Range("N34").Value = -3464.2
$ws.Range("H39").Value = 90917390
$ws.Range("I39").Value = 111117810
$ws.Range("J39").Value = 15500
$ws.Range("K39").Value = 111117810
$ws.Range("L39").Value = 15500
$ws.Range("M39").Value = -111117419
$ws.Range("N39").Value = -16282
$ws.Range("H49").Value = 90917390
$ws.Range("I49").Value = 111117810
$ws.Range("J49").Value = 15500
$ws.Range("K49").Value = 111117810
$ws.Range("L49").Value = 15500
$ws.Range("M49").Value = -111117628
$ws.Range("N49").Value = -15864
$ws.Range("H64").Value = 11200
$ws.Range("J64").Value = 11200
$ws.Range("L64").Value = 11200
$ws.Range("N64").Value = -11696
$ws.Range("H67").Value = 11200
$ws.Range("J67").Value = 11200
$ws.Range("L67").Value = 11200
$ws.Range("N67").Value = -12916
$ws.Range("H107").Value = 705.7857
$ws.Range("I107").Value = 712.2857
$ws.Range("K107").Value = 712.2857
$ws.Range("M107").Value = 1207.7143
$ws.Range("H122").Value = 1039.75
$ws.Range("I122").Value = 1039.75
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3119.25
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -669.25
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 78.71429000000001
$ws.Range("I7").Value = 78.71429000000001
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 236.14287
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -124.14287
$ws.Range("N7").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 3058.75
$ws.Range("I33").Value = 3000
$ws.Range("J33").Value = 3078.3333
$ws.Range("K33").Value = 3000
$ws.Range("L33").Value = 3078.3333
$ws.Range("M33").Value = -2748
$ws.Range("N33").Value = -3582.3333
$ws.Range("H52").Value = 50000
$ws.Range("J52").Value = 50000
$ws.Range("L52").Value = 50000
$ws.Range("N52").Value = -50518
$ws.Range("H123").Value = 14034.637
$ws.Range("J123").Value = 14034.637
$ws.Range("L123").Value = 14034.637
$ws.Range("N123").Value = -18934.637
$ws.Range("H134").Value = 41325.5
$ws.Range("J134").Value = 41325.5
$ws.Range("L134").Value = 123976.5
$ws.Range("N134").Value = -129046.5
$ws.Range("H136").Value = 13383
$ws.Range("J136").Value = 13383
$ws.Range("L136").Value = 40149
$ws.Range("N136").Value = -45249

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1268.16
$ws.Range("I68").Value = 958.1053000000001
$ws.Range("J68").Value = 2250
$ws.Range("K68").Value = 958.1053000000001
$ws.Range("L68").Value = 2250
$ws.Range("M68").Value = -209.1053000000001
$ws.Range("N68").Value = -3748
$ws.Range("H71").Value = 1268.16
$ws.Range("I71").Value = 958.1053000000001
$ws.Range("J71").Value = 2250
$ws.Range("K71").Value = 4790.5265
$ws.Range("L71").Value = 11250
$ws.Range("M71").Value = -1046.5265
$ws.Range("N71").Value = -18738
$ws.Range("H122").Value = 8337862.5
$ws.Range("J122").Value = 3157.5
$ws.Range("L122").Value = 9472.5
$ws.Range("N122").Value = -14372.5
$ws.Range("H132").Value = 4539.6
$ws.Range("I132").Value = 5168.731
$ws.Range("J132").Value = 2722.111
$ws.Range("K132").Value = 15506.193
$ws.Range("L132").Value = 8166.333
$ws.Range("M132").Value = -12976.193
$ws.Range("N132").Value = -13226.333
$ws.Range("H135").Value = 66214.5
$ws.Range("J135").Value = 66214.5
$ws.Range("L135").Value = 66214.5
$ws.Range("N135").Value = -76354.5
$ws.Range("H136").Value = 4667
$ws.Range("I136").Value = 5423.5
$ws.Range("J136").Value = 1857.1428
$ws.Range("K136").Value = 16270.5
$ws.Range("L136").Value = 5571.428400000001
$ws.Range("M136").Value = -13720.5
$ws.Range("N136").Value = -10671.4284

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2232.9412
$ws.Range("I132").Value = 1613.05
$ws.Range("K132").Value = 4839.15
$ws.Range("M132").Value = -2309.15
